$wb = $excel.ActiveWorkbook
$wsTable = $wb.Worksheets.Item("Table")
$wsGraph = $wb.Worksheets.Item("Graph")

# --- Append a new day of data (row 125, date 2023-01-18 / serial 44917) ---
# Mirrors the pattern of the preceding rows (121-124): literal inputs in
# A/B/D/E/F/G/I/J/K, SUM formulas in C and H.
$wsTable.Cells.Item(125, 1).Value = 44917
$wsTable.Cells.Item(125, 2).Value = 9
$wsTable.Cells.Item(125, 3).Formula = "=SUM(D125:F125)"
$wsTable.Cells.Item(125, 4).Value = 1105
$wsTable.Cells.Item(125, 5).Value = 371
$wsTable.Cells.Item(125, 6).Value = 197
$wsTable.Cells.Item(125, 7).Value = 40
$wsTable.Cells.Item(125, 8).Formula = "=SUM(I125:K125)"
# Keep the formula cell's style consistent with the rest of column H (no
# explicit style index), instead of whatever format got auto-applied.
$wsTable.Cells.Item(125, 8).Style = "Normal"
$wsTable.Cells.Item(125, 9).Value = 3281
$wsTable.Cells.Item(125, 10).Value = 434
$wsTable.Cells.Item(125, 11).Value = 385

# --- Correct the "Swift code" day count for 2023-01-11 .. 2023-01-16 (rows 121-124) ---
$wsTable.Cells.Item(121, 7).Value = 40
$wsTable.Cells.Item(122, 7).Value = 40
$wsTable.Cells.Item(123, 7).Value = 40
$wsTable.Cells.Item(124, 7).Value = 40

# --- Move the active selection on the Table sheet to P118 ---
[void]$wsTable.Range("P118").Select()

# --- Make "Table" the active/selected tab instead of "Graph" ---
[void]$wsTable.Activate()
